$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Asset Class"
$ws.Range("B1").Value = "Nama"
$ws.Range("A2").Value = "Infrastructure"
$ws.Range("B2").Value = "Container & Ramps"

$ws.Range("B6").Select()
